$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that must be forced to text (values that look numeric) to avoid Excel
# auto-converting the string into a number/date, matching the original inlineStr cells.
$textCells = @{
    'D2' = '29.870.36'
    'D3' = '1.900.07'
    'D4' = '0.9999'
    'D5' = '0.7619'
    'D9' = '25.59'
    'D10' = '0.06842'
    'D11' = '0.07962'
    'D12' = '1.905.41'
    'D13' = '0.7436'
    'D14' = '5.172'
    'D15' = '90.92'
    'D16' = '29.866.85'
    'D17' = '13.97'
    'D18' = '5.955'
    'D19' = '243.19'
    'D20' = '0.000007690'
    'D21' = '0.9998'
    'D22' = '1.000'
    'D23' = '6.949'
    'D24' = '166.03'
    'D25' = '9.229'
    'D27' = '0.1290'
    'D28' = '2.046'
    'D29' = '1.415'
    'D30' = '1.513'
    'D31' = '4.257'
    'D32' = '4.073'
    'D33' = '0.05231'
    'D34' = '1.255'
    'D35' = '0.7277'
    'D36' = '2.714'
    'D38' = '2.781'
    'D39' = '6.163'
    'D40' = '0.4409'
    'D41' = '72.00'
    'D42' = '0.9998'
    'D43' = '1.886'
    'D44' = '0.8266'
    'D45' = '7.644'
    'D46' = '100.01'
    'D47' = '9.750'
    'D48' = '2.055.61'
    'D51' = '1.473'
}

foreach ($ref in $textCells.Keys) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $textCells[$ref]
    $c.Style = "Normal"
}

# Remaining cells (coin name / link / volume percentage) are safe to set directly,
# they are not ambiguous with Excel numeric/date auto-detection.
$plainCells = @{
    'E2' = '  -0.16%  '
    'E3' = '  +0.53%  '
    'E4' = '  -0.39%  '
    'E5' = '  +5.11%  '
    'E6' = '  -0.73%  '
    'E7' = '  -0.27%  '
    'E8' = '  -0.84%  '
    'E9' = '  -1.66%  '
    'E10' = '  -0.42%  '
    'E11' = '  +0.24%  '
    'E12' = '  +0.83%  '
    'E13' = '  -2.86%  '
    'E14' = '  -1.10%  '
    'E15' = '  +0.06%  '
    'E16' = '  -0.22%  '
    'E17' = '  -0.64%  '
    'E18' = '  +3.93%  '
    'E19' = '  +2.10%  '
    'E20' = '  -0.62%  '
    'E21' = '  -0.30%  '
    'E22' = '  -0.37%  '
    'E23' = '  +2.00%  '
    'E24' = '  +0.67%  '
    'E25' = '  -0.29%  '
    'E26' = '  -0.96%  '
    'E27' = '  +1.62%  '
    'E28' = '  +2.04%  '
    'E29' = '  +4.25%  '
    'E30' = '  -1.11%  '
    'E31' = '  -0.62%  '
    'E32' = '  +0.51%  '
    'E33' = '  +3.66%  '
    'E34' = '  -0.80%  '
    'E36' = '  -0.73%  '
    'E37' = '  +0.64%  '
    'E38' = '  +0.42%  '
    'E39' = '  -2.30%  '
    'E41' = '  -2.96%  '
    'E42' = '  -0.06%  '
    'E43' = '  -1.42%  '
    'E44' = '  -1.00%  '
    'E45' = '  +1.09%  '
    'B46' = 'Quant'
    'C46' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E46' = '  -0.58%  '
    'B47' = 'EnergySwap'
    'C47' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E47' = '  +0.49%  '
    'E48' = '  +0.93%  '
    'E49' = '  -3.47%  '
    'E50' = '  -0.14%  '
    'E51' = '  +2.13%  '
}

foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}
